$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 75
$ws1.Range("F3").Value = 402
$ws1.Range("F4").Value = 449
$ws1.Range("F6").Value = 18
$ws1.Range("F8").Value = 14141
$ws1.Range("F9").Value = 128
$ws1.Range("F10").Value = 103
$ws1.Range("F11").Value = 5677
$ws1.Range("F12").Value = 581
$ws1.Range("F15").Value = 54
$ws1.Range("F16").Value = 1228
$ws1.Range("F17").Value = 3
$ws1.Range("F19").Value = 167
$ws1.Range("F20").Value = 769
$ws1.Range("F21").Value = 2915
$ws1.Range("F22").Value = 49
$ws1.Range("F23").Value = 10471
$ws1.Range("F24").Value = 1191
$ws1.Range("F25").Value = 46
$ws1.Range("F26").Value = 66
$ws1.Range("F27").Value = 3714
$ws1.Range("F29").Value = 64

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 75
$ws4.Range("F3").Value = 402
$ws4.Range("F5").Value = 449
$ws4.Range("F7").Value = 18
$ws4.Range("F9").Value = 14141
$ws4.Range("F10").Value = 128
$ws4.Range("F11").Value = 103
$ws4.Range("F12").Value = 5677
$ws4.Range("F13").Value = 581
$ws4.Range("F16").Value = 54
$ws4.Range("F17").Value = 1228
$ws4.Range("F18").Value = 3
$ws4.Range("F20").Value = 167
$ws4.Range("F21").Value = 769
$ws4.Range("F22").Value = 2915
$ws4.Range("F23").Value = 49
$ws4.Range("F25").Value = 10471
$ws4.Range("F26").Value = 1191
$ws4.Range("F27").Value = 46
$ws4.Range("F28").Value = 66
$ws4.Range("F29").Value = 3714
$ws4.Range("F31").Value = 64
